# M16 slide: combine the two bulleted lists (tasks + tools) into the
# single explanatory sentence used after the #937/#1019 rework.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Remember the shape's current (un-autofitted) height so we can restore
# it after the text edit - PowerPoint's spAutoFit recalculates box height
# as paragraphs are removed, but the original box size is unrelated to
# this content change.
$originalHeight = $shape.Height

# Drop every bullet paragraph except the first one, merging all the
# "taken"/"tools" bullets away.
for ($i = $tr.Paragraphs().Count; $i -ge 2; $i--) {
    $tr.Paragraphs($i, 1).Delete()
}
# The engine leaves a trailing empty paragraph shell behind after the
# loop above (Paragraphs().Count under-reports it) - remove it too.
$tr.Paragraphs(2, 1).Delete()

# Replace the remaining (first) paragraph's text with the new combined
# explanation, editing the existing run in place.
$tr.Paragraphs(1, 1).Text = "Voor vastgestelde taken bij het ontwikkelen, onderhouden en operationeel beheren van software, stelt ICTU het gebruik van tools verplicht. ICTU adviseert per taak specifieke tools en ondersteunt projecten bij het gebruik daarvan."

# Restore the shape's original height (undo the autofit shrink).
$shape.Height = $originalHeight
